$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.945.10"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "2.240.84"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'268.27"
$ws.Range("E5").Value = "  +4.87%  "
$ws.Range("D6").Value = "'88.21"
$ws.Range("E6").Value = "  +13.81%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").Value = "'46.13"
$ws.Range("E10").Value = "  +8.72%  "
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "'7.59"
$ws.Range("E12").Value = "  +8.62%  "
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").Value = "2.564.55"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "'14.83"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "2.241.64"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "'0.797"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "43.857.70"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "'6.04"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'70.13"
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("D23").Value = "'232.89"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  -4.16%  "
$ws.Range("D25").Value = "'2.60"
$ws.Range("E25").Value = "  +18.58%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "'3.54"
$ws.Range("E28").Value = "  +6.02%  "
$ws.Range("D29").Value = "'40.36"
$ws.Range("E29").Value = "  -5.31%  "
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").Value = "'175.43"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "'0.0913"
$ws.Range("E32").Value = "  +5.33%  "
$ws.Range("D33").Value = "'20.69"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("D37").Value = "'0.0358"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "'3.34"
$ws.Range("E39").Value = "  +17.01%  "
$ws.Range("D40").Value = "'12.56"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").Value = "'65.48"
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("D42").Value = "'2.12"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").Value = "'0.0996"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("D46").Value = "'100.31"
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").Value = "'8.34"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "'1.22"
$ws.Range("E48").Value = "  +7.47%  "
$ws.Range("D49").Value = "'1.14"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.54"
$ws.Range("E50").Value = "  +4.90%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.440"
$ws.Range("E51").Value = "  -9.99%  "
